# edit.ps1
#
# Applies the edit described by the commit:
#   "Almost finished the proposal.  Summary needs to be extended and
#    Unique Selling Position needs finishing."
#
# Two changes to docs/proposal.docx:
#   1) The "Game Play:" section paragraph is rewritten with a longer
#      description of the quest system, and picks up the same
#      double-spaced / contextual-spacing paragraph formatting used by
#      the other body paragraphs.
#   2) The "Summary:" section, previously three empty paragraphs, gets
#      its first paragraph filled in with closing summary text (two of
#      the blank paragraphs are consumed in the process).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: find the (1-based) paragraph index whose text starts with a
# given needle. Using text search rather than a hard-coded index keeps
# this script robust to paragraph-count drift between edits.
# ---------------------------------------------------------------------
function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) "Game Play:" section paragraph - replace the two terse runs with
#    the new, longer description, and add the spacing/contextualSpacing
#    paragraph properties.
# ---------------------------------------------------------------------
$gamePlayIdx = Find-ParagraphIndex("Game play will subsist")
$gamePlayPara = $d.Paragraphs($gamePlayIdx)

$gamePlayXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:contextualSpacing/></w:pPr><w:r><w:t xml:space="preserve">The player will be free to move throughout the office. He will complete quests to advance in the game, which will be accessed by talking to NPCs in the office. Quests will consist of </w:t></w:r><w:r><w:t>delivering</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:t>retrieving</w:t></w:r><w:r><w:t xml:space="preserve"> items or relaying messages between NPCs.</w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r></w:p>
"@

[void]$gamePlayPara.Range.InsertXML($gamePlayXml)

# ---------------------------------------------------------------------
# 2) "Summary:" section - the first of the trailing blank paragraphs
#    gets the new summary text...
# ---------------------------------------------------------------------
$summaryHeadingIdx = Find-ParagraphIndex("Summary:")
$summaryParaIdx = $summaryHeadingIdx + 1
$summaryPara = $d.Paragraphs($summaryParaIdx)

$summaryXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>In summary, Intern's Quest will be a witty adventure RPG that will draw players in with it</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> dialog and charm.</w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:t xml:space="preserve">Our team has the knowledge, experience, and drive to develop Intern's Quest. </w:t></w:r></w:p>
"@

[void]$summaryPara.Range.InsertXML($summaryXml)

# ...and two of the now-redundant trailing blank paragraphs are removed
# (five blank paragraphs following "Summary:" become one filled-in
# paragraph plus two blanks). Paragraphs are deleted one at a time,
# each as its own single-paragraph range, since deleting a single range
# spanning exactly two empty paragraphs only consumes one of them.
for ($n = 1; $n -le 2; $n++) {
    $target = $d.Paragraphs($summaryParaIdx + 1)
    $delRange = $d.Range($target.Range.Start, $target.Range.End)
    [void]$delRange.Delete()
}

Write-Output "Edit complete."
